$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$gValues = @{
    2 = 1
    3 = 3
    4 = 0
    5 = 0
    6 = 1
    7 = 2
    8 = 0
    9 = 4
    10 = 2
    11 = 0
    12 = 0
    13 = 3
    14 = 2
    15 = 3
    16 = 3
    17 = 4
    18 = 5
    19 = 3
    20 = 2
    21 = 2
    22 = 1
    24 = 0
    25 = 2
    26 = 1
    27 = 2
    28 = 1
    29 = 1
    30 = 1
    31 = 0
    32 = 2
    33 = 0
    34 = 3
    35 = 0
    36 = 0
    37 = 3
    38 = 2
    39 = 2
    40 = 1
    41 = 0
    42 = 1
    43 = 3
    44 = 3
    45 = 2
    46 = 2
    47 = 2
    48 = 0
    49 = 1
    50 = 1
    51 = 0
    52 = 1
    53 = 1
    54 = 1
    55 = 1
    56 = 3
    57 = 1
    58 = 3
    59 = 1
    60 = 3
    61 = 2
    62 = 2
    63 = 2
    64 = 4
    65 = 3
    66 = 2
    67 = 3
    68 = 1
    69 = 1
    70 = 1
    71 = 3
    72 = 1
    73 = 1
    74 = 0
    75 = 0
    76 = 5
    77 = 2
    78 = 2
    79 = 1
    80 = 0
}

foreach ($row in $gValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $gValues[$row]
}
